# Auto-generated Word COM-interop script to update multiplication answers
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "188×7=1316" -> "823×4=3292"
$cell = $t.Cell(1, 1)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("188×7=1316")) {
    Write-Host "WARNING: Row 1 Col 1 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "823×4=3292"

# Row 1, Col 2: "691×5=3455" -> "627×8=5016"
$cell = $t.Cell(1, 2)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("691×5=3455")) {
    Write-Host "WARNING: Row 1 Col 2 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "627×8=5016"

# Row 1, Col 3: "614×6=3684" -> "796×5=3980"
$cell = $t.Cell(1, 3)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("614×6=3684")) {
    Write-Host "WARNING: Row 1 Col 3 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "796×5=3980"

# Row 1, Col 4: "560×5=2800" -> "285×3=855"
$cell = $t.Cell(1, 4)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("560×5=2800")) {
    Write-Host "WARNING: Row 1 Col 4 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "285×3=855"

# Row 1, Col 5: "867×6=5202" -> "998×7=6986"
$cell = $t.Cell(1, 5)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("867×6=5202")) {
    Write-Host "WARNING: Row 1 Col 5 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "998×7=6986"

# Row 5, Col 1: "233×2=466" -> "212×3=636"
$cell = $t.Cell(5, 1)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("233×2=466")) {
    Write-Host "WARNING: Row 5 Col 1 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "212×3=636"

# Row 5, Col 2: "595×4=2380" -> "972×5=4860"
$cell = $t.Cell(5, 2)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("595×4=2380")) {
    Write-Host "WARNING: Row 5 Col 2 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "972×5=4860"

# Row 5, Col 3: "141×5=705" -> "855×8=6840"
$cell = $t.Cell(5, 3)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("141×5=705")) {
    Write-Host "WARNING: Row 5 Col 3 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "855×8=6840"

# Row 5, Col 4: "519×6=3114" -> "911×5=4555"
$cell = $t.Cell(5, 4)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("519×6=3114")) {
    Write-Host "WARNING: Row 5 Col 4 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "911×5=4555"

# Row 5, Col 5: "548×5=2740" -> "271×8=2168"
$cell = $t.Cell(5, 5)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("548×5=2740")) {
    Write-Host "WARNING: Row 5 Col 5 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "271×8=2168"

# Row 10, Col 1: "181×4=724" -> "733×7=5131"
$cell = $t.Cell(10, 1)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("181×4=724")) {
    Write-Host "WARNING: Row 10 Col 1 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "733×7=5131"

# Row 10, Col 2: "823×6=4938" -> "339×8=2712"
$cell = $t.Cell(10, 2)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("823×6=4938")) {
    Write-Host "WARNING: Row 10 Col 2 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "339×8=2712"

# Row 10, Col 3: "885×9=7965" -> "560×8=4480"
$cell = $t.Cell(10, 3)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("885×9=7965")) {
    Write-Host "WARNING: Row 10 Col 3 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "560×8=4480"

# Row 10, Col 4: "718×5=3590" -> "179×3=537"
$cell = $t.Cell(10, 4)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("718×5=3590")) {
    Write-Host "WARNING: Row 10 Col 4 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "179×3=537"

# Row 10, Col 5: "596×3=1788" -> "725×6=4350"
$cell = $t.Cell(10, 5)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("596×3=1788")) {
    Write-Host "WARNING: Row 10 Col 5 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "725×6=4350"

# Row 15, Col 1: "543×9=4887" -> "462×4=1848"
$cell = $t.Cell(15, 1)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("543×9=4887")) {
    Write-Host "WARNING: Row 15 Col 1 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "462×4=1848"

# Row 15, Col 2: "137×3=411" -> "323×4=1292"
$cell = $t.Cell(15, 2)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("137×3=411")) {
    Write-Host "WARNING: Row 15 Col 2 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "323×4=1292"

# Row 15, Col 3: "489×9=4401" -> "361×9=3249"
$cell = $t.Cell(15, 3)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("489×9=4401")) {
    Write-Host "WARNING: Row 15 Col 3 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "361×9=3249"

# Row 15, Col 4: "632×7=4424" -> "417×5=2085"
$cell = $t.Cell(15, 4)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("632×7=4424")) {
    Write-Host "WARNING: Row 15 Col 4 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "417×5=2085"

# Row 15, Col 5: "845×3=2535" -> "363×5=1815"
$cell = $t.Cell(15, 5)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("845×3=2535")) {
    Write-Host "WARNING: Row 15 Col 5 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "363×5=1815"

# Row 20, Col 1: "522×8=4176" -> "852×8=6816"
$cell = $t.Cell(20, 1)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("522×8=4176")) {
    Write-Host "WARNING: Row 20 Col 1 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "852×8=6816"

# Row 20, Col 2: "179×3=537" -> "144×7=1008"
$cell = $t.Cell(20, 2)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("179×3=537")) {
    Write-Host "WARNING: Row 20 Col 2 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "144×7=1008"

# Row 20, Col 3: "894×7=6258" -> "819×6=4914"
$cell = $t.Cell(20, 3)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("894×7=6258")) {
    Write-Host "WARNING: Row 20 Col 3 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "819×6=4914"

# Row 20, Col 4: "141×4=564" -> "912×8=7296"
$cell = $t.Cell(20, 4)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("141×4=564")) {
    Write-Host "WARNING: Row 20 Col 4 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "912×8=7296"

# Row 20, Col 5: "385×9=3465" -> "243×2=486"
$cell = $t.Cell(20, 5)
$cellText = $cell.Range.Text
if ($cellText -notmatch [regex]::Escape("385×9=3465")) {
    Write-Host "WARNING: Row 20 Col 5 did not contain expected text. Found:" $cellText
}
$cell.Range.Text = "243×2=486"

Write-Host "All replacements applied."
